$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 183.6
$ws.Range("I6").Value = 79.5
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 238.5
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -126.5
$ws.Range("N6").Value = -2024
$ws.Range("H62").Value = 11371553
$ws.Range("I62").Value = 16676748
$ws.Range("J62").Value = 3278.4285
$ws.Range("K62").Value = 16676748
$ws.Range("L62").Value = 3278.4285
$ws.Range("M62").Value = -16676124
$ws.Range("N62").Value = -4526.4285
$ws.Range("H65").Value = 11371553
$ws.Range("I65").Value = 16676748
$ws.Range("J65").Value = 3278.4285
$ws.Range("K65").Value = 83383740
$ws.Range("L65").Value = 16392.1425
$ws.Range("M65").Value = -83380620
$ws.Range("N65").Value = -22632.1425

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1904.2174
$ws.Range("I3").Value = 279.4
$ws.Range("J3").Value = 2355.5557
$ws.Range("K3").Value = 279.4
$ws.Range("L3").Value = 2355.5557
$ws.Range("M3").Value = -164.4
$ws.Range("N3").Value = -2585.5557
$ws.Range("H32").Value = 1535.41
$ws.Range("I32").Value = 1489.236
$ws.Range("J32").Value = 1909
$ws.Range("K32").Value = 1489.236
$ws.Range("L32").Value = 1909
$ws.Range("M32").Value = -1202.236
$ws.Range("N32").Value = -2483
$ws.Range("H61").Value = 1742.5416
$ws.Range("I61").Value = 1379.4375
$ws.Range("J61").Value = 2468.75
$ws.Range("K61").Value = 1379.4375
$ws.Range("L61").Value = 2468.75
$ws.Range("M61").Value = -1167.4375
$ws.Range("N61").Value = -2892.75
$ws.Range("H74").Value = 1456.5135
$ws.Range("I74").Value = 1478.5454
$ws.Range("J74").Value = 1274.75
$ws.Range("K74").Value = 1478.5454
$ws.Range("L74").Value = 1274.75
$ws.Range("M74").Value = -604.5454
$ws.Range("N74").Value = -3022.75
$ws.Range("H77").Value = 1456.5135
$ws.Range("I77").Value = 1478.5454
$ws.Range("J77").Value = 1274.75
$ws.Range("K77").Value = 7392.727
$ws.Range("L77").Value = 6373.75
$ws.Range("M77").Value = -3024.727
$ws.Range("N77").Value = -15109.75
$ws.Range("H132").Value = 18577190
$ws.Range("I132").Value = 19608674
$ws.Range("J132").Value = 9809588
$ws.Range("K132").Value = 58826022
$ws.Range("L132").Value = 29428764
$ws.Range("M132").Value = -58823492
$ws.Range("N132").Value = -29433824
$ws.Range("H136").Value = 1742.5416
$ws.Range("I136").Value = 1379.4375
$ws.Range("J136").Value = 2468.75
$ws.Range("K136").Value = 4138.3125
$ws.Range("L136").Value = 7406.25
$ws.Range("M136").Value = -1588.3125
$ws.Range("N136").Value = -12506.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1738493
$ws.Range("I134").Value = 710.3674
$ws.Range("J134").Value = 7415249.5
$ws.Range("K134").Value = 2131.1022
$ws.Range("L134").Value = 22245748.5
$ws.Range("M134").Value = 403.8978000000002
$ws.Range("N134").Value = -22250818.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H31").Value = 1373.6986
$ws.Range("I31").Value = 1110.0646
$ws.Range("J31").Value = 1568.2858
$ws.Range("K31").Value = 1110.0646
$ws.Range("L31").Value = 1568.2858
$ws.Range("M31").Value = -815.0645999999999
$ws.Range("N31").Value = -2158.2858
$ws.Range("H34").Value = 1373.6986
$ws.Range("I34").Value = 1110.0646
$ws.Range("J34").Value = 1568.2858
$ws.Range("K34").Value = 1110.0646
$ws.Range("L34").Value = 1568.2858
$ws.Range("M34").Value = -908.0645999999999
$ws.Range("N34").Value = -1972.2858
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H58").Value = 12500691
$ws.Range("I58").Value = 22727856
$ws.Range("J58").Value = 823.75
$ws.Range("K58").Value = 22727856
$ws.Range("L58").Value = 823.75
$ws.Range("M58").Value = -22727653
$ws.Range("N58").Value = -1229.75
$ws.Range("H134").Value = 806.93335
$ws.Range("I134").Value = 688.34485
$ws.Range("J134").Value = 1211.5294
$ws.Range("K134").Value = 2065.03455
$ws.Range("L134").Value = 3634.5882
$ws.Range("M134").Value = 469.9654500000001
$ws.Range("N134").Value = -8704.5882
$ws.Range("H136").Value = 12500691
$ws.Range("I136").Value = 22727856
$ws.Range("J136").Value = 823.75
$ws.Range("K136").Value = 68183568
$ws.Range("L136").Value = 2471.25
$ws.Range("M136").Value = -68181018
$ws.Range("N136").Value = -7571.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 27779330
$ws.Range("I137").Value = 31250956
$ws.Range("J137").Value = 23811760
$ws.Range("K137").Value = 93752868
$ws.Range("L137").Value = 71435280
$ws.Range("M137").Value = -93747768
$ws.Range("N137").Value = -71445480

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4957.227
$ws.Range("I70").Value = 4088.6667
$ws.Range("J70").Value = 6818.4287
$ws.Range("K70").Value = 4088.6667
$ws.Range("L70").Value = 6818.4287
$ws.Range("M70").Value = -3818.6667
$ws.Range("N70").Value = -7358.4287
$ws.Range("H73").Value = 4957.227
$ws.Range("I73").Value = 4088.6667
$ws.Range("J73").Value = 6818.4287
$ws.Range("K73").Value = 4088.6667
$ws.Range("L73").Value = 6818.4287
$ws.Range("M73").Value = -3152.6667
$ws.Range("N73").Value = -8690.4287
$ws.Range("H132").Value = 4890.1567
$ws.Range("I132").Value = 3151.4
$ws.Range("J132").Value = 11212.909
$ws.Range("K132").Value = 9454.200000000001
$ws.Range("L132").Value = 33638.727
$ws.Range("M132").Value = -6924.200000000001
$ws.Range("N132").Value = -38698.727

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5544.8965
$ws.Range("I122").Value = 9192.462
$ws.Range("J122").Value = 2581.25
$ws.Range("K122").Value = 27577.386
$ws.Range("L122").Value = 7743.75
$ws.Range("M122").Value = -25127.386
$ws.Range("N122").Value = -12643.75
$ws.Range("H136").Value = 39410930
$ws.Range("I136").Value = 7145388
$ws.Range("J136").Value = 111112136
$ws.Range("K136").Value = 21436164
$ws.Range("L136").Value = 333336408
$ws.Range("M136").Value = -21433614
$ws.Range("N136").Value = -333341508

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 38462012
$ws.Range("I113").Value = 47619436
$ws.Range("J113").Value = 838.6
$ws.Range("K113").Value = 142858308
$ws.Range("L113").Value = 2515.8
$ws.Range("M113").Value = -142856138
$ws.Range("N113").Value = -6855.8
$ws.Range("H132").Value = 6344895.5
$ws.Range("I132").Value = 17945.066
$ws.Range("J132").Value = 27786228
$ws.Range("K132").Value = 53835.198
$ws.Range("L132").Value = 83358684
$ws.Range("M132").Value = -51305.198
$ws.Range("N132").Value = -83363744
$ws.Range("H136").Value = 13164192
$ws.Range("I136").Value = 16673441
$ws.Range("J136").Value = 4508.125
$ws.Range("K136").Value = 50020323
$ws.Range("L136").Value = 13524.375
$ws.Range("M136").Value = -18624.375
$ws.Range("N136").Value = -18624.375
